$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '30.006.32'
$ws.Range("E2").Value = '  +1.37%  '

# Row 3
$ws.Range("D3").Value = '1.637.40'
$ws.Range("E3").Value = '  +2.35%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.11%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '215.64'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.57%  '

# Row 6
$ws.Range("E6").Value = '  +1.59%  '

# Row 7
$ws.Range("E7").Value = '  -0.13%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '29.75'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +10.78%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.262'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +4.27%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0615'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +2.42%  '

# Row 11
$ws.Range("E11").Value = '  +0.48%  '

# Row 12
$ws.Range("D12").Value = '1.869.67'
$ws.Range("E12").Value = '  +2.23%  '

# Row 13
$ws.Range("D13").Value = '1.635.42'
$ws.Range("E13").Value = '  +2.20%  '

# Row 14
$ws.Range("B14").Value = 'Polygon'
$ws.Range("C14").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.579'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +7.29%  '

# Row 15
$ws.Range("B15").Value = 'Chainlink'
$ws.Range("C15").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '9.63'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +26.22%  '

# Row 16
$ws.Range("E16").Value = '  +4.60%  '

# Row 17
$ws.Range("D17").Value = '30.033.71'
$ws.Range("E17").Value = '  +1.44%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '64.93'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.90%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '249.24'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +3.13%  '

# Row 20
$ws.Range("E20").Value = '  +2.34%  '

# Row 21
$ws.Range("E21").Value = '  -0.11%  '

# Row 22
$ws.Range("E22").Value = '  +5.42%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.77'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +5.68%  '

# Row 24
$ws.Range("E24").Value = '  +1.71%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '159.37'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +2.76%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '15.75'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +2.53%  '

# Row 27
$ws.Range("E27").Value = '  +2.83%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.68'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +4.37%  '

# Row 29
$ws.Range("E29").Value = '  -0.12%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0492'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.90%  '

# Row 31
$ws.Range("E31").Value = '  +6.50%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.40'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +5.48%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.22'
$ws.Range("D33").Style = "Normal"

# Row 34
$ws.Range("D34").Value = '1.435.79'
$ws.Range("E34").Value = '  +0.87%  '

# Row 35
$ws.Range("E35").Value = '  +7.97%  '

# Row 36
$ws.Range("E36").Value = '  +1.68%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.86'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.43%  '

# Row 38
$ws.Range("E38").Value = '  +2.22%  '

# Row 39
$ws.Range("E39").Value = '  -0.30%  '

# Row 40
$ws.Range("B40").Value = 'Aave'
$ws.Range("C40").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '77.32'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +16.59%  '

# Row 41
$ws.Range("B41").Value = 'ImmutableX'
$ws.Range("C41").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.560'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.90%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.840'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +3.82%  '

# Row 43
$ws.Range("E43").Value = '  +1.47%  '

# Row 44
$ws.Range("B44").Value = 'Kaspa'
$ws.Range("C44").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0499'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.95%  '

# Row 45
$ws.Range("B45").Value = 'BitcoinSV'
$ws.Range("C45").Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '55.04'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.62%  '

# Row 46
$ws.Range("E46").Value = '  +4.63%  '

# Row 48
$ws.Range("E48").Value = '  +2.17%  '

# Row 49
$ws.Range("D49").Value = '1.776.83'
$ws.Range("E49").Value = '  +2.12%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '90.48'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +4.93%  '

# Row 51
$ws.Range("E51").Value = '  +5.18%  '
